$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last populated row in the sheet is row 23. Two new rows of exercise-log
# data are appended below it (rows 24 and 25). Start by cloning the
# formatting (date number format on column A, "Dialog" font styling on the
# numeric columns, etc.) from row 23 so the new rows look consistent with
# the rest of the table.
$ws.Range("A23:I23").Copy() | Out-Null
$ws.Range("A24:I25").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 24 — 2020-06-27
$ws.Cells.Item(24, 1).Value = 44009
$ws.Cells.Item(24, 2).Value = 98.5
$ws.Cells.Item(24, 3).Value = 104
$ws.Cells.Item(24, 4).Value = 0.9
$ws.Cells.Item(24, 5).Value = "IMPROVED"
$ws.Cells.Item(24, 6).Value = 0.1
$ws.Cells.Item(24, 7).Value = 83.3
$ws.Cells.Item(24, 8).Value = 29.2
$ws.Cells.Item(24, 9).Value = "OVERWEIGHT"

# Row 25 — 2020-06-28
$ws.Cells.Item(25, 1).Value = 44010
$ws.Cells.Item(25, 2).Value = 98
$ws.Cells.Item(25, 3).Value = 104
$ws.Cells.Item(25, 4).Value = 0.9
$ws.Cells.Item(25, 5).Value = "SAME"
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 83.8
$ws.Cells.Item(25, 8).Value = 29.3
$ws.Cells.Item(25, 9).Value = "OVERWEIGHT"
